$wb = $excel.ActiveWorkbook

# Rename the three "Tab_IND_*" sheets to "Tab_IN_*"
$wb.Worksheets.Item("Tab_IND_5-10years").Name = "Tab_IN_5-10years"
$wb.Worksheets.Item("Tab_IND_10-15years").Name = "Tab_IN_10-15years"
$wb.Worksheets.Item("Tab_IND_15+years").Name = "Tab_IN_15+years"

# Sheet "GenAI Startup Overview 2024": zoom 87 -> 125
$ws1 = $wb.Worksheets.Item("GenAI Startup Overview 2024")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 125

# Sheet "Competitor Insights": selection D21 (sqref A1:XFD1048576) -> D7:H8
$ws3 = $wb.Worksheets.Item("Competitor Insights")
$ws3.Activate()
$ws3.Range("D7:H8").Select()

# Sheet "Tab_EMEA_5-10years": selection H42 -> C2 (no longer the tab-selected sheet)
$ws4 = $wb.Worksheets.Item("Tab_EMEA_5-10years")
$ws4.Activate()
$ws4.Range("C2").Select()

# Sheet "Tab_IN_15+years" (last sheet, formerly "Tab_IND_15+years"): becomes the
# active/selected tab; selection E11 (sqref A2:E11) -> Z37
$ws12 = $wb.Worksheets.Item("Tab_IN_15+years")
$ws12.Activate()
$ws12.Range("Z37").Select()
